# Update FuelPrices at 2025-04-07 08:29
# Column A (prices) and Column B (dates) swap places, and a new data row
# (r=19) is appended with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row: swap "Date" and "MLBSO00" between A1 and B1 ----
$ws.Range("A1").Value = "MLBSO00"
$ws.Range("B1").Value = "Date"

# ---- Data rows 2-18: swap date/price values between columns A and B ----
# Date values (were in column A, move to column B)
$dates = @(45734, 45733, 45730, 45729, 45728, 45735, 45736, 45737, 45740, 45741, 45742, 45734, 45743, 45744, 45748, 45749, 45750)
# Price values (were in column B, move to column A)
$prices = @(806.651, 810.465, 810.465, 810.465, 810.465, 806.651, 806.651, 806.651, 806.651, 806.651, 806.651, 806.651, 806.651, 800.9299999999999, 800.9299999999999, 800.9299999999999, 797.116)

for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 1).Value = $prices[$i]
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $dates[$i]
}

# ---- New row 19: append the latest reading ----
$ws.Cells.Item(19, 1).ClearFormats()
$ws.Cells.Item(19, 1).Value = 791.395
$ws.Cells.Item(19, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(19, 2).Value = 45751
$ws.Cells.Item(19, 3).Value = 730.444

# Row 18's date cell now uses the regular (non-last-row) date format,
# since row 19 is the new last row.
$ws.Cells.Item(18, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
